# Update the cryptocurrency price/volume table on Sheet1 to the latest
# scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '27.064.65'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''215.44'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +1.81%  '
$ws.Range('E9').Value = '  +5.40%  '
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').Value = '1.912.30'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '1.676.38'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').Value = '''65.97'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').Value = '27.055.55'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '''237.57'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('D20').Value = '0.0₃0741'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').Value = '''9.34'
$ws.Range('E23').Value = '  +1.82%  '
$ws.Range('E24').Value = '  -1.96%  '
$ws.Range('D25').Value = '''146.71'
$ws.Range('E25').Value = '  +0.55%  '
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('D27').Value = '''16.35'
$ws.Range('E27').Value = '  +2.61%  '
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = '''0.0497'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''3.35'
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.546.72'
$ws.Range('E33').Value = '  +5.95%  '
$ws.Range('E34').Value = '  +1.70%  '
$ws.Range('E35').Value = '  +2.72%  '
$ws.Range('D36').Value = '''0.599'
$ws.Range('E36').Value = '  +3.23%  '
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').Value = '''0.923'
$ws.Range('E38').Value = '  +2.75%  '
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('E40').Value = '  +1.86%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('D42').Value = '''67.62'
$ws.Range('D44').Value = '''2.25'
$ws.Range('E44').Value = '  -1.95%  '
$ws.Range('D45').Value = '1.821.35'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('D46').Value = '''0.782'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').Value = '''90.67'
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('E49').Value = '  +1.97%  '
$ws.Range('E50').Value = '  +2.65%  '
$ws.Range('D51').Value = '''8.07'
$ws.Range('E51').Value = '  +5.32%  '
